# Domain_Model.xlsx - 1st User story element complete
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename "TakeAway" user story to "Takeaway"
$ws.Range("A2").Value = "Takeaway"

# Replace the completed methods list with the first implemented method
$ws.Range("B2").Value = "print_menu; "

# Row height was only tall to fit the old, longer text - let it auto-fit again
$ws.Rows.Item(2).AutoFit() | Out-Null

# Move / restore the active selection as left by the author
$ws.Range("D18").Select() | Out-Null
